{"js": "const body = context.document.body;\nconst ooxml = \"<?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:pPr><w:rPr><w:lang w:val=\\\"ru-RU\\\"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\nconst r = body.insertOoxml(ooxml, Word.InsertLocation.end);\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nreturn \"count:\" + paras.items.length;\n", "ps1": "$d = $word.ActiveDocument\n$e = $d.Content.End\n$r = $d.Range($e, $e)\n$r.InsertParagraphAfter()\n$d.Paragraphs.Count\n"}
